$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifting existing rows down)
$ws.Rows.Item(2).Insert()

# Populate the new row with the new election entry
$ws.Range("A2").Value = "Presidencia Municipal 15"
$ws.Range("B2").Value = "pm_15"
$ws.Range("C2").Value = "#669bbc"

# Update the active selection to match the target state
$ws.Range("A2:C2").Select()
